$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.512.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.635.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("E12").Value = '  +3.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.089.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +13.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.507.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.630.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '347.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0797'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("E30").Value = '  +5.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  +11.19%  '
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.88%  '
$ws.Range("E37").Value = '  +4.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '330.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.38%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.58%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.00%  '
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("E42").Value = '  +4.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '133.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.58%  '
$ws.Range("E48").Value = '  +2.60%  '
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("E51").Value = '  +0.93%  '
